$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: donnemartin/data-science-ipython-notebooks
$ws.Range("A3").Value = "donnemartin"
$ws.Range("B3").Value = "data-science-ipython-notebooks"
$ws.Range("C3").Value = 42027.81839120371
$ws.Range("D3").Value = 44504.84664351852
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 49025
$ws.Range("G3").Value = "Python"
$ws.Range("H3").Value = 4
$ws.Range("I3").Value = 12
$ws.Range("J3").Value = 1627
$ws.Range("K3").Value = 21825
$ws.Range("L3").Value = 6809
$ws.Range("M3").Value = 31
$ws.Range("N3").Value = 543
$ws.Range("O3").Value = 53
$ws.Range("P3").Value = 2
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = "Data science Python notebooks: Deep learning (TensorFlow, Theano, Caffe, Keras), scikit-learn, Kaggle, big data (Spark, Hadoop MapReduce, HDFS), matplotlib, pandas, NumPy, SciPy, Python essentials, AWS, and various command lines."

# Row 4: virgili0/Virgilio
$ws.Range("A4").Value = "virgili0"
$ws.Range("B4").Value = "Virgilio"
$ws.Range("C4").Value = 43536.78166666667
$ws.Range("D4").Value = 44488.26559027778
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 35112
$ws.Range("G4").Value = "Jupyter Notebook"
$ws.Range("H4").Value = 6
$ws.Range("I4").Value = 31
$ws.Range("J4").Value = 785
$ws.Range("K4").Value = 13189
$ws.Range("L4").Value = 2509
$ws.Range("M4").Value = 37
$ws.Range("N4").Value = 1363
$ws.Range("O4").Value = 123
$ws.Range("P4").Value = 12
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = "Your new Mentor for Data Science E-Learning."

# Row 5: microsoft/Data-Science-For-Beginners
$ws.Range("A5").Value = "microsoft"
$ws.Range("B5").Value = "Data-Science-For-Beginners"
$ws.Range("C5").Value = 44258.68767361111
$ws.Range("D5").Value = 44519.76876157407
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 62297
$ws.Range("G5").Value = "Jupyter Notebook"
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 87
$ws.Range("J5").Value = 144
$ws.Range("K5").Value = 6991
$ws.Range("L5").Value = 967
$ws.Range("M5").Value = 72
$ws.Range("N5").Value = 924
$ws.Range("O5").Value = 215
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = "10 Weeks, 20 Lessons, Data Science for All!"

# Row 6: joelgrus/data-science-from-scratch
$ws.Range("A6").Value = "joelgrus"
$ws.Range("B6").Value = "data-science-from-scratch"
$ws.Range("C6").Value = 41952.10513888889
$ws.Range("D6").Value = 44313.96119212963
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 769
$ws.Range("G6").Value = "Python"
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 13
$ws.Range("J6").Value = 624
$ws.Range("K6").Value = 6499
$ws.Range("L6").Value = 3682
$ws.Range("M6").Value = 70
$ws.Range("N6").Value = 89
$ws.Range("O6").Value = 39
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = "code for Data Science From Scratch book"

# Row 7: firmai/industry-machine-learning
$ws.Range("A7").Value = "firmai"
$ws.Range("B7").Value = "industry-machine-learning"
$ws.Range("C7").Value = 43588.2140625
$ws.Range("D7").Value = 44489.48476851852
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 3067
$ws.Range("G7").Value = "Jupyter Notebook"
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 385
$ws.Range("K7").Value = 6065
$ws.Range("L7").Value = 992
$ws.Range("M7").Value = 2
$ws.Range("N7").Value = 157
$ws.Range("O7").Value = 5
$ws.Range("P7").Value = 1
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = "A curated list of applied machine learning and data science notebooks and libraries across different industries (by @firmai)"

# Row 8: fengdu78/Data-Science-Notes
$ws.Range("A8").Value = "fengdu78"
$ws.Range("B8").Value = "Data-Science-Notes"
$ws.Range("C8").Value = 43675.16336805555
$ws.Range("D8").Value = 44424.47130787037
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 52628
$ws.Range("G8").Value = "Jupyter Notebook"
$ws.Range("H8").Value = 2
$ws.Range("I8").Value = 1
$ws.Range("J8").Value = 220
$ws.Range("K8").Value = 5960
$ws.Range("L8").Value = 2595
$ws.Range("M8").Value = 13
$ws.Range("N8").Value = 70
$ws.Range("O8").Value = 5
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 0
$ws.Range("S8").Value = "数据科学的笔记以及资料搜集"

# Row 9: rushter/data-science-blogs
$ws.Range("A9").Value = "rushter"
$ws.Range("B9").Value = "data-science-blogs"
$ws.Range("C9").Value = 42221.49877314815
$ws.Range("D9").Value = 44411.99362268519
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 348
$ws.Range("G9").Value = "Python"
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 86
$ws.Range("J9").Value = 470
$ws.Range("K9").Value = 5723
$ws.Range("L9").Value = 1718
$ws.Range("M9").Value = 7
$ws.Range("N9").Value = 359
$ws.Range("O9").Value = 114
$ws.Range("P9").Value = 1
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 0
$ws.Range("S9").Value = "A curated list of data science blogs"

# Row 10: drivendata/cookiecutter-data-science
$ws.Range("A10").Value = "drivendata"
$ws.Range("B10").Value = "cookiecutter-data-science"
$ws.Range("C10").Value = 42307.80552083333
$ws.Range("D10").Value = 44508.65208333333
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 724
$ws.Range("G10").Value = "Python"
$ws.Range("H10").Value = 4
$ws.Range("I10").Value = 45
$ws.Range("J10").Value = 98
$ws.Range("K10").Value = 5185
$ws.Range("L10").Value = 1704
$ws.Range("M10").Value = 124
$ws.Range("N10").Value = 145
$ws.Range("O10").Value = 134
$ws.Range("P10").Value = 7
$ws.Range("Q10").Value = 1
$ws.Range("R10").Value = 1
$ws.Range("S10").Value = "A logical, reasonably standardized, but flexible project structure for doing and sharing data science work."

# Row 11: Netflix/metaflow
$ws.Range("A11").Value = "Netflix"
$ws.Range("B11").Value = "metaflow"
$ws.Range("C11").Value = 43725.74195601852
$ws.Range("D11").Value = 44519.87226851852
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 9725
$ws.Range("G11").Value = "Python"
$ws.Range("H11").Value = 5
$ws.Range("I11").Value = 463
$ws.Range("J11").Value = 227
$ws.Range("K11").Value = 4968
$ws.Range("L11").Value = 433
$ws.Range("M11").Value = 337
$ws.Range("N11").Value = 385
$ws.Range("O11").Value = 491
$ws.Range("P11").Value = 33
$ws.Range("Q11").Value = 36
$ws.Range("R11").Value = 33
$ws.Range("S11").Value = ":rocket: Build and manage real-life data science projects with ease!"
